$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03)
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
